$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 3: новая строка — тот же пункт отгрузки/код, но дата на день позже
# и объём фасовки "округлён"/разнесён (было 8 -> 5 + 3)
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = (Get-Date -Year 2022 -Month 6 -Day 12).Date
$ws.Range("D3").Value = 5

# Row 4: новая строка — ещё на день позже
$ws.Range("A4").Value = $ws.Range("A2").Value2
$ws.Range("B4").Value = $ws.Range("B2").Value2
$ws.Range("C4").Value = (Get-Date -Year 2022 -Month 6 -Day 13).Date
$ws.Range("D4").Value = 3
